$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix calculation of change: "Initial biomass" value is reset to 0 ---
$ws.Range("B8").Value = 0

# --- Add "Final total biomass (tonnes CO2e)" header in C12, matching B12's style ---
$ws.Range("B12").Copy()
$ws.Range("C12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C12").Value = "Final total biomass`n(tonnes CO2e)"

# Make the trailing "2" subscript (and keep it bold, matching the sibling header in B12)
$sub = $ws.Range("C12").Characters(31, 1)
$sub.Font.Subscript = $true
$sub.Font.Bold = $true

# Keep the closing "e)" bold (non-subscript) as its own run
$tail = $ws.Range("C12").Characters(32, 2)
$tail.Font.Bold = $true

# --- New "Final total biomass" value cell C13, matching B13's style, value 0 ---
$ws.Range("B13").Copy()
$ws.Range("C13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C13").Value = 0

# --- Update the active selection to B9 ---
[void]$ws.Range("B9").Select()

$excel.CutCopyMode = $false
